# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value for each affected sheet.
$updates = @{
    2  = 146
    3  = 1676
    4  = 772
    6  = 30
    7  = 11787
    11 = 397
    13 = 839
    14 = 13432
    15 = 13313
    20 = 266
    21 = 91
    23 = 152
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
